# Apply updated "dSF" (column F) values for a set of rows in Sheet1.
# These correspond to a re-pull of data where column F values were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    3  = -5
    9  = -2
    12 = -6
    13 = -7
    14 = -6
    16 = 0
    17 = 0
    18 = 5
    19 = -2
    20 = 1
    22 = -4
    27 = 3
    30 = -1
    33 = 2
    34 = -5
    39 = -2
    44 = 1
    45 = -8
    47 = -1
    49 = 2
    51 = -3
    55 = -3
    58 = -2
    62 = -2
    63 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
